# This script rotates the B..G column values among groups of duplicate
# product rows (rows that share the same item description but represent
# different stock "lots"). Within each group, row i receives the
# B..G values that previously belonged to row i+1 (wrapping around to the
# first row of the group). Column A (serial no.) and H..M are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row groups (1-based worksheet rows) that get their data cyclically
# rotated, taken from the diff between the original and edited workbook.
$groups = @(
    @(149, 150),
    @(161, 162, 163),
    @(264, 265),
    @(316, 317, 318),
    @(346, 347),
    @(351, 352),
    @(355, 356),
    @(372, 373),
    @(375, 376),
    @(431, 432),
    @(579, 580),
    @(583, 584),
    @(586, 587),
    @(720, 721),
    @(859, 860),
    @(889, 890)
)

# Columns (by index, B=2 .. G=7) whose values are rotated.
$cols = 2..7

foreach ($group in $groups) {
    $n = $group.Length

    # Snapshot current values for every row/column in this group first,
    # so that writing new values doesn't clobber data we still need to read.
    $snapshot = @{}
    foreach ($row in $group) {
        $rowValues = @{}
        foreach ($col in $cols) {
            $rowValues[$col] = $ws.Cells.Item($row, $col).Value2
        }
        $snapshot[$row] = $rowValues
    }

    for ($i = 0; $i -lt $n; $i++) {
        $destRow = $group[$i]
        $srcRow = $group[($i + 1) % $n]
        $srcValues = $snapshot[$srcRow]
        foreach ($col in $cols) {
            $ws.Cells.Item($destRow, $col).Value = $srcValues[$col]
        }
    }
}
